# Update cryptocurrency price/volume data in the worksheet (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "51.847.44"
$ws.Range("E2").Value = "  +1.28%  "
# Row 3: Ethereum
$ws.Range("D3").Value = "3.010.62"
$ws.Range("E3").Value = "  +3.29%  "
# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.01%  "
# Row 5: BNB
$ws.Range("D5").Value = "'381.73"
$ws.Range("E5").Value = "  +4.90%  "
# Row 6: Solana
$ws.Range("D6").Value = "'107.14"
$ws.Range("E6").Value = "  +2.38%  "
# Row 7: XRP
$ws.Range("E7").Value = "  +1.46%  "
# Row 8: USDC
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.04%  "
# Row 9: Cardano
$ws.Range("E9").Value = "  +2.21%  "
# Row 10: Avalanche
$ws.Range("E10").Value = "  +2.36%  "
# Row 11: TRON
$ws.Range("E11").Value = "  +0.70%  "
# Row 13: Chainlink
$ws.Range("D13").Value = "'18.81"
$ws.Range("E13").Value = "  +1.91%  "
# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.485.58"
$ws.Range("E14").Value = "  +3.48%  "
# Row 15: Polkadot
$ws.Range("D15").Value = "'7.57"
$ws.Range("E15").Value = "  +3.12%  "
# Row 16: WrappedEther
$ws.Range("D16").Value = "3.004.56"
$ws.Range("E16").Value = "  +3.33%  "
# Row 17: Polygon
$ws.Range("E17").Value = "  +2.40%  "
# Row 18: WrappedBTC
$ws.Range("D18").Value = "51.879.47"
$ws.Range("E18").Value = "  +1.48%  "
# Row 19: ImmutableX
$ws.Range("D19").Value = "'3.39"
$ws.Range("E19").Value = "  +2.43%  "
# Row 20: Uniswap
$ws.Range("E20").Value = "  +3.33%  "
# Row 21: InternetComputer(DFINITY)
$ws.Range("D21").Value = "'13.18"
$ws.Range("E21").Value = "  +1.36%  "
# Row 22: ShibaInu
$ws.Range("E22").Value = "  +1.84%  "
# Row 23: Litecoin
$ws.Range("D23").Value = "'69.09"
$ws.Range("E23").Value = "  +1.31%  "
# Row 24: BitcoinCash
$ws.Range("D24").Value = "'265.02"
$ws.Range("E24").Value = "  +2.21%  "
# Row 25: PancakeSwap
$ws.Range("E25").Value = "  +4.84%  "
# Row 26: Kaspa
$ws.Range("E26").Value = "  +0.35%  "
# Row 27: RenderToken
$ws.Range("D27").Value = "'7.24"
$ws.Range("E27").Value = "  +17.49%  "
# Row 28: Filecoin
$ws.Range("E28").Value = "  +4.79%  "
# Row 29: EthereumClassic
$ws.Range("D29").Value = "'26.27"
$ws.Range("E29").Value = "  +1.45%  "
# Row 30: Dai
$ws.Range("E30").Value = "  -0.12%  "
# Row 31: Hedera
$ws.Range("E31").Value = "  +1.62%  "
# Row 32: Cosmos
$ws.Range("E32").Value = "  +0.80%  "
# Row 33: InjectiveProtocol
$ws.Range("D33").Value = "'35.20"
$ws.Range("E33").Value = "  +0.68%  "
# Row 34: OKB
$ws.Range("D34").Value = "'51.57"
$ws.Range("E34").Value = "  +1.71%  "
# Row 35: Toncoin
$ws.Range("E35").Value = "  -2.80%  "
# Row 36: VeChain
$ws.Range("D36").Value = "'0.0458"
$ws.Range("E36").Value = "  +8.70%  "
# Row 37: FirstDigitalUSD
$ws.Range("E37").Value = "  +0.12%  "
# Row 38: LidoDAOToken
$ws.Range("D38").Value = "'3.16"
$ws.Range("E38").Value = "  +0.64%  "
# Row 39: Celestia
$ws.Range("E39").Value = "  +4.37%  "
# Row 40: Stacks
$ws.Range("E40").Value = "  -5.41%  "
# Row 41: ARBITRUM
$ws.Range("E41").Value = "  +0.98%  "
# Row 42: Stellar
$ws.Range("E42").Value = "  +2.93%  "
# Row 43: Monero
$ws.Range("D43").Value = "'124.59"
$ws.Range("E43").Value = "  +4.57%  "
# Row 44: EnergySwap
$ws.Range("D44").Value = "'22.50"
$ws.Range("E44").Value = "  +0.68%  "
# Row 45: WEMIXToken
$ws.Range("E45").Value = "  -1.61%  "
# Row 46: TheGraph
$ws.Range("E46").Value = "  +18.43%  "
# Row 47: ApeXProtocol
$ws.Range("D47").Value = "'2.39"
$ws.Range("E47").Value = "  +5.43%  "
# Row 48: Maker
$ws.Range("D48").Value = "2.062.44"
$ws.Range("E48").Value = "  -0.01%  "
# Row 49: NEARProtocol
$ws.Range("E49").Value = "  +3.65%  "
# Row 50: BEAM
$ws.Range("E50").Value = "  +14.48%  "
# Row 51: FraxShare
$ws.Range("D51").Value = "'8.88"
$ws.Range("E51").Value = "  +3.44%  "
